$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing "extr" rows (8-15) before we overwrite anything, then
# rewrite them two rows further down (10-17) to make room for two new
# contingency lines ("line7" / "line8") inserted right after "line6".
$orig = @{}
for ($i = 8; $i -le 15; $i++) {
    $orig[$i] = @($ws.Range("A$i").Value2, $ws.Range("B$i").Value2, $ws.Range("C$i").Value2, $ws.Range("D$i").Value2, $ws.Range("E$i").Value2)
}
for ($i = 15; $i -ge 8; $i--) {
    $dst = $i + 2
    $vals = $orig[$i]
    $ws.Range("A$dst").Value = $vals[0] + 2
    $ws.Range("B$dst").Value = $vals[1]
    $ws.Range("C$dst").Value = $vals[2]
    $ws.Range("D$dst").Value = $vals[3]
    $ws.Range("E$dst").Value = $vals[4]
}

# The two new rows at the bottom (16, 17) need the same style as column A's
# other data cells (bold/centered/bordered), since they are brand new rows.
$ws.Range("A2").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)

# New row 8: line7
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# New row 9: line8
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# extr7 (now row 16) flips its in_service flag from false to true
$ws.Range("E16").Value = $true
